$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new text value.
# All target values are written as literal text (matching the source
# data, which stores prices/percentages/labels as text, not numbers).
$updates = @(
    @{r=2; c=4; v='67.733.48'},
    @{r=2; c=5; v='  +8.97%  '},
    @{r=3; c=4; v='3.522.64'},
    @{r=3; c=5; v='  +10.89%  '},
    @{r=4; c=5; v='  -0.17%  '},
    @{r=5; c=4; v='190.68'},
    @{r=5; c=5; v='  +13.44%  '},
    @{r=6; c=4; v='556.74'},
    @{r=6; c=5; v='  +10.02%  '},
    @{r=7; c=4; v='3.519.04'},
    @{r=7; c=5; v='  +10.79%  '},
    @{r=8; c=5; v='  +5.05%  '},
    @{r=9; c=5; v='  -0.08%  '},
    @{r=10; c=4; v='0.640'},
    @{r=10; c=5; v='  +9.60%  '},
    @{r=11; c=4; v='57.05'},
    @{r=11; c=5; v='  +6.38%  '},
    @{r=12; c=4; v='0.152'},
    @{r=12; c=5; v='  +18.98%  '},
    @{r=13; c=4; v='0.0000277'},
    @{r=13; c=5; v='  +12.03%  '},
    @{r=14; c=5; v='  +9.41%  '},
    @{r=15; c=4; v='4.067.59'},
    @{r=15; c=5; v='  +9.87%  '},
    @{r=16; c=4; v='3.518.13'},
    @{r=16; c=5; v='  +10.29%  '},
    @{r=17; c=4; v='68.478.41'},
    @{r=17; c=5; v='  +10.14%  '},
    @{r=18; c=5; v='  +8.70%  '},
    @{r=19; c=4; v='18.38'},
    @{r=19; c=5; v='  +9.89%  '},
    @{r=20; c=5; v='  +12.43%  '},
    @{r=21; c=5; v='  +9.65%  '},
    @{r=22; c=4; v='409.44'},
    @{r=22; c=5; v='  +15.16%  '},
    @{r=23; c=4; v='3.97'},
    @{r=23; c=5; v='  +10.15%  '},
    @{r=24; c=4; v='11.76'},
    @{r=24; c=5; v='  +11.98%  '},
    @{r=25; c=4; v='84.88'},
    @{r=25; c=5; v='  +9.11%  '},
    @{r=26; c=4; v='4.23'},
    @{r=26; c=5; v='  +12.55%  '},
    @{r=27; c=2; v='ImmutableX'},
    @{r=27; c=3; v='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{r=27; c=4; v='2.91'},
    @{r=27; c=5; v='  +13.64%  '},
    @{r=28; c=2; v='LEO'},
    @{r=28; c=3; v='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{r=28; c=4; v='6.16'},
    @{r=28; c=5; v='  +1.64%  '},
    @{r=29; c=4; v='11.91'},
    @{r=29; c=5; v='  +9.77%  '},
    @{r=30; c=4; v='8.64'},
    @{r=30; c=5; v='  +8.79%  '},
    @{r=31; c=4; v='30.60'},
    @{r=31; c=5; v='  +11.23%  '},
    @{r=32; c=4; v='686.72'},
    @{r=32; c=5; v='  +12.64%  '},
    @{r=33; c=4; v='6.87'},
    @{r=33; c=5; v='  +9.01%  '},
    @{r=34; c=4; v='11.75'},
    @{r=34; c=5; v='  +8.46%  '},
    @{r=35; c=5; v='  +10.66%  '},
    @{r=36; c=4; v='60.29'},
    @{r=36; c=5; v='  +7.37%  '},
    @{r=37; c=2; v='PEPE'},
    @{r=37; c=3; v='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'},
    @{r=37; c=4; v='0.0₃0836'},
    @{r=37; c=5; v='  +28.00%  '},
    @{r=38; c=2; v='InjectiveProtocol'},
    @{r=38; c=3; v='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'},
    @{r=38; c=4; v='39.10'},
    @{r=38; c=5; v='  +10.48%  '},
    @{r=39; c=4; v='0.405'},
    @{r=39; c=5; v='  +10.23%  '},
    @{r=40; c=5; v='  +0.26%  '},
    @{r=41; c=4; v='3.43'},
    @{r=41; c=5; v='  +28.48%  '},
    @{r=42; c=5; v='  +13.37%  '},
    @{r=43; c=4; v='2.77'},
    @{r=43; c=5; v='  +18.30%  '},
    @{r=44; c=4; v='3.03'},
    @{r=44; c=5; v='  +18.96%  '},
    @{r=45; c=2; v='Maker'},
    @{r=45; c=3; v='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{r=45; c=4; v='3.071.87'},
    @{r=45; c=5; v='  +9.23%  '},
    @{r=46; c=2; v='FirstDigitalUSD'},
    @{r=46; c=3; v='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'},
    @{r=46; c=4; v='0.997'},
    @{r=46; c=5; v='  -0.16%  '},
    @{r=47; c=5; v='  +11.84%  '},
    @{r=48; c=4; v='9.29'},
    @{r=48; c=5; v='  +25.61%  '},
    @{r=49; c=4; v='2.76'},
    @{r=49; c=5; v='  +7.75%  '},
    @{r=50; c=4; v='3.23'},
    @{r=50; c=5; v='  +11.15%  '},
    @{r=51; c=4; v='0.131'},
    @{r=51; c=5; v='  +8.70%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.r, $u.c)
    # Force text storage even for values that look numeric (e.g. "190.68"),
    # then restore the default "Normal" style so no stray number format
    # is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value2 = $u.v
    $cell.Style = "Normal"
}
